# Fruta / hortaliza, semanal
# Insert 2 new weekly price rows (Brócoli, Terminal Hortofrutícola Agro Chillán)
# right above the former last-but-one data block, pushing the existing rows
# (and the former last row) down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 270-271; everything from row 270 down (incl. the
# former row 296) shifts down to 272-298. Formatting (incl. the date style
# on column D) is inherited from the row above, matching the surrounding
# records.
$ws.Range("A270:A271").EntireRow.Insert()

# --- New row 270 ---
$ws.Range("A270").Value = 7
$ws.Range("B270").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C270").Value = "Ñuble"
$ws.Range("D270").Value = 44783
$ws.Range("E270").Value = 16
$ws.Range("F270").Value = 100112023
$ws.Range("G270").Value = "Brócoli"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 200
$ws.Range("K270").Value = 800
$ws.Range("L270").Value = 900
$ws.Range("M270").Value = 850
$ws.Range("N270").Value = "$/unidad"
$ws.Range("O270").Value = "Provincia de Diguillín"
$ws.Range("P270").Value = 850
$ws.Range("Q270").Value = 1
$ws.Range("R270").Value = "Hortaliza"

# --- New row 271 ---
$ws.Range("A271").Value = 7
$ws.Range("B271").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C271").Value = "Ñuble"
$ws.Range("D271").Value = 44783
$ws.Range("E271").Value = 16
$ws.Range("F271").Value = 100112023
$ws.Range("G271").Value = "Brócoli"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Segunda"
$ws.Range("J271").Value = 200
$ws.Range("K271").Value = 700
$ws.Range("L271").Value = 700
$ws.Range("M271").Value = 700
$ws.Range("N271").Value = "$/unidad"
$ws.Range("O271").Value = "Provincia de Diguillín"
$ws.Range("P271").Value = 700
$ws.Range("Q271").Value = 1
$ws.Range("R271").Value = "Hortaliza"
